$d = $word.ActiveDocument

# --- Step 1: remove the "Meta description" paragraph (the one whose text
# starts with "Meta description") -------------------------------------------
$metaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Meta description*") {
        $metaIndex = $i
        break
    }
}
if ($metaIndex -gt 0) {
    $d.Paragraphs.Item($metaIndex).Range.Delete()
}

# --- Step 2: insert a new bold paragraph ("Play Gangsterz Slot Machine for
# Free - Review") right before the final paragraph (the italic image-prompt
# paragraph that starts with "Please create a feature image") --------------
$promptIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Please create a feature image*") {
        $promptIndex = $i
        break
    }
}
if ($promptIndex -gt 0) {
    $promptPara = $d.Paragraphs.Item($promptIndex)
    $promptPara.Range.InsertParagraphBefore()

    $newPara = $d.Paragraphs.Item($promptIndex)
    $newRange = $newPara.Range

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Gangsterz Slot Machine for Free - Review</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $newRange.InsertXML($xml)

    # --- Step 3: replace the text of the (now shifted) prompt paragraph with
    # the new description copy, keeping its italic formatting and leading
    # empty run intact -------------------------------------------------------
    $promptIndex = $promptIndex + 1
    $promptPara = $d.Paragraphs.Item($promptIndex)
    $pr = $promptPara.Range
    $bodyRange = $d.Range($pr.Start, $pr.End - 1)
    $bodyRange.Text = "Explore Gangsterz slot machine game with winning features and high RTP value. Play free to discover the game's design and strategies."
}

Write-Host "Edit complete. Paragraph count:" $d.Paragraphs.Count
